# EPBDS-4790 [Web Services] Incorrect result of running test table in case
# project has depended module
#
# The test table on Sheet1 gains a "properties / version / 0.0.2" row after
# every test-case header row, so that the dependent-module test results
# reflect the module version the test was produced against.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new (blank) row right after each existing test-case header row.
# These row numbers are expressed in terms of the *current* sheet layout at
# the moment each Insert() runs, i.e. each subsequent insertion point already
# accounts for the rows inserted before it.
$insertAt = @(4, 9, 14, 19, 23)
foreach ($r in $insertAt) {
    $ws.Rows.Item($r).Insert()
}

# Fill in the new "properties" header rows with the property name/value pairs.
$propertyRows = @(4, 9, 14, 19, 23)
foreach ($r in $propertyRows) {
    $cellRange = $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 4))
    $cellRange.Borders.LineStyle = 1
    $cellRange.HorizontalAlignment = -4108

    $ws.Cells.Item($r, 2).Value = "properties"
    $ws.Cells.Item($r, 3).Value = "version"
    $ws.Cells.Item($r, 4).Value = "0.0.2"
}

# Match the new selection left behind by the edit.
$ws.Range("B23:D23").Select()

